# Updates the cryptos list (Price / Volume(1h) columns, and a couple of
# rows whose rank changed so Coin/Link/Price/Volume all moved).
# Values are written with a leading "'" so Excel stores them as literal
# text (matching the source data's inline-string cells, e.g. "42.788.70")
# instead of auto-converting number-looking strings into numeric values.
# Resetting .Style back to "Normal" afterwards drops the quote-prefix
# formatting so the cell keeps the sheet's default (unstyled) look.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.788.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.47%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.561.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.38%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'302.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'96.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +3.99%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.575"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.56%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.09%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.544"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.18%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'36.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.66%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0809"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.79%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +9.95%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'7.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.22%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.573.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.12%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.880"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.92%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'14.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.98%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'42.861.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.61%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'13.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +5.26%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +2.20%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.72%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'71.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'257.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.28%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'2.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.51%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -1.01%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'28.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -4.31%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.00%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'39.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +9.19%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'10.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.54%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -1.56%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'6.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.09%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'156.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +3.91%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +1.07%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.08%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'26.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +11.23%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.0798"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.85%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -2.34%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +1.73%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'18.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +15.06%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.41%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +1.90%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +28.92%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'VeChain"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.0307"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.18%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'NEARProtocol"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'3.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.24%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'FirstDigitalUSD"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.06%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'Maker"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'2.061.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.36%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D47").Value = "'9.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +6.41%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'76.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +11.37%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.811.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.58%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'103.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.80%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +3.07%  "
$ws.Range("E51").Style = "Normal"
